$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1465.3552
$ws.Range("I15").Value = 1465.3552
$ws.Range("K15").Value = 4396.0656
$ws.Range("M15").Value = -4227.0656
$ws.Range("H80").Value = 1949.9166
$ws.Range("I80").Value = 2504.3333
$ws.Range("J80").Value = 286.66666
$ws.Range("K80").Value = 7512.999899999999
$ws.Range("L80").Value = 859.9999799999999
$ws.Range("M80").Value = -6514.999899999999
$ws.Range("N80").Value = -2855.99998
$ws.Range("H83").Value = 1949.9166
$ws.Range("I83").Value = 2504.3333
$ws.Range("J83").Value = 286.66666
$ws.Range("K83").Value = 22538.9997
$ws.Range("L83").Value = 2579.99994
$ws.Range("M83").Value = -17546.9997
$ws.Range("N83").Value = -12563.99994
$ws.Range("H125").Value = 1803.2
$ws.Range("J125").Value = 1803.2
$ws.Range("L125").Value = 16228.8
$ws.Range("N125").Value = -21148.8
$ws.Range("H132").Value = 4009.95
$ws.Range("I132").Value = 3344.5
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 10033.5
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -7503.5
$ws.Range("N132").Value = -35057
$ws.Range("H138").Value = 1919.6897
$ws.Range("I138").Value = 1432.7812
$ws.Range("K138").Value = 4298.3436
$ws.Range("M138").Value = 841.6563999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1644.3334
$ws.Range("I2").Value = 1573.2
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 1573.2
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -1460.2
$ws.Range("N2").Value = -2226
$ws.Range("H32").Value = 532093.5600000001
$ws.Range("I32").Value = 609729.7
$ws.Range("J32").Value = 23145.777
$ws.Range("K32").Value = 609729.7
$ws.Range("L32").Value = 23145.777
$ws.Range("M32").Value = -609442.7
$ws.Range("N32").Value = -23719.777
$ws.Range("H45").Value = 4505.615
$ws.Range("I45").Value = 4725
$ws.Range("J45").Value = 4154.6
$ws.Range("K45").Value = 4725
$ws.Range("L45").Value = 4154.6
$ws.Range("M45").Value = -4348
$ws.Range("N45").Value = -4908.6
$ws.Range("H63").Value = 5684.769
$ws.Range("J63").Value = 6111.5
$ws.Range("L63").Value = 6111.5
$ws.Range("N63").Value = -7483.5
$ws.Range("H66").Value = 5684.769
$ws.Range("J66").Value = 6111.5
$ws.Range("L66").Value = 30557.5
$ws.Range("N66").Value = -37421.5
$ws.Range("H116").Value = 1644.3334
$ws.Range("I116").Value = 1573.2
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 1573.2
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 720.8
$ws.Range("N116").Value = -6588
$ws.Range("H122").Value = 48737.43
$ws.Range("I122").Value = 53551.895
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 160655.685
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -158205.685
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1644.3334
$ws.Range("I3").Value = 1573.2
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 1573.2
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -1459.2
$ws.Range("N3").Value = -2228
$ws.Range("H123").Value = 28221.334
$ws.Range("J123").Value = 28221.334
$ws.Range("L123").Value = 28221.334
$ws.Range("N123").Value = -38021.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1986.5714
$ws.Range("I99").Value = 1849.5
$ws.Range("J99").Value = 2041.4
$ws.Range("K99").Value = 1849.5
$ws.Range("L99").Value = 2041.4
$ws.Range("M99").Value = -351.5
$ws.Range("N99").Value = -5037.4
$ws.Range("H122").Value = 1904.2307
$ws.Range("I122").Value = 1844.375
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5533.125
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3083.125
$ws.Range("N122").Value = -10900
$ws.Range("H126").Value = 1986.5714
$ws.Range("I126").Value = 1849.5
$ws.Range("J126").Value = 2041.4
$ws.Range("K126").Value = 5548.5
$ws.Range("L126").Value = 6124.200000000001
$ws.Range("M126").Value = -3078.5
$ws.Range("N126").Value = -11064.2
$ws.Range("H132").Value = 7578164
$ws.Range("I132").Value = 2622.5
$ws.Range("J132").Value = 13891115
$ws.Range("K132").Value = 7867.5
$ws.Range("L132").Value = 41673345
$ws.Range("M132").Value = -5337.5
$ws.Range("N132").Value = -41678405

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 646
$ws.Range("I5").Value = 470.5
$ws.Range("J5").Value = 1207.6
$ws.Range("K5").Value = 1411.5
$ws.Range("L5").Value = 3622.8
$ws.Range("M5").Value = -1299.5
$ws.Range("N5").Value = -3846.8
$ws.Range("H23").Value = 78947500
$ws.Range("I23").Value = 129.66667
$ws.Range("K23").Value = 389.00001
$ws.Range("M23").Value = -154.00001
$ws.Range("H107").Value = 568.4545000000001
$ws.Range("I107").Value = 393.2857
$ws.Range("J107").Value = 875
$ws.Range("K107").Value = 1179.8571
$ws.Range("L107").Value = 2625
$ws.Range("M107").Value = 740.1428999999998
$ws.Range("N107").Value = -6465
$ws.Range("H122").Value = 5858.6
$ws.Range("I122").Value = 714.8570999999999
$ws.Range("J122").Value = 8628.308000000001
$ws.Range("K122").Value = 6433.7139
$ws.Range("L122").Value = 77654.77200000001
$ws.Range("M122").Value = -3983.7139
$ws.Range("N122").Value = -82554.77200000001
$ws.Range("H132").Value = 2243.7942
$ws.Range("J132").Value = 2409.9583
$ws.Range("L132").Value = 21689.6247
$ws.Range("N132").Value = -26749.6247
$ws.Range("H135").Value = 646
$ws.Range("I135").Value = 470.5
$ws.Range("J135").Value = 1207.6
$ws.Range("K135").Value = 4234.5
$ws.Range("L135").Value = 10868.4
$ws.Range("M135").Value = -1699.5
$ws.Range("N135").Value = -15938.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 251753.25
$ws.Range("I113").Value = 501000
$ws.Range("J113").Value = 2506.5
$ws.Range("K113").Value = 501000
$ws.Range("L113").Value = 2506.5
$ws.Range("M113").Value = -498830
$ws.Range("N113").Value = -6846.5
$ws.Range("H122").Value = 35715750
$ws.Range("I122").Value = 1586.25
$ws.Range("J122").Value = 250000750
$ws.Range("K122").Value = 4758.75
$ws.Range("L122").Value = 750002250
$ws.Range("M122").Value = -2308.75
$ws.Range("N122").Value = -750007150
$ws.Range("H132").Value = 3170.074
$ws.Range("I132").Value = 3064.3076
$ws.Range("J132").Value = 3268.2856
$ws.Range("K132").Value = 9192.9228
$ws.Range("L132").Value = 9804.856800000001
$ws.Range("M132").Value = -6662.9228
$ws.Range("N132").Value = -14864.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 99715
$ws.Range("J36").Value = 99715
$ws.Range("L36").Value = 99715
$ws.Range("N36").Value = -100839
$ws.Range("H46").Value = 907.1429000000001
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 936.36365
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 936.36365
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -1312.36365
$ws.Range("H68").Value = 1636.5
$ws.Range("I68").Value = 1525
$ws.Range("J68").Value = 2751.5
$ws.Range("K68").Value = 1525
$ws.Range("L68").Value = 2751.5
$ws.Range("M68").Value = -776
$ws.Range("N68").Value = -4249.5
$ws.Range("H71").Value = 1636.5
$ws.Range("I71").Value = 1525
$ws.Range("J71").Value = 2751.5
$ws.Range("K71").Value = 7625
$ws.Range("L71").Value = 13757.5
$ws.Range("M71").Value = -3881
$ws.Range("N71").Value = -21245.5
$ws.Range("H122").Value = 3307.6191
$ws.Range("I122").Value = 3136.3635
$ws.Range("J122").Value = 3496
$ws.Range("K122").Value = 9409.0905
$ws.Range("L122").Value = 10488
$ws.Range("M122").Value = -6959.0905
$ws.Range("N122").Value = -15388
$ws.Range("H136").Value = 16668918
$ws.Range("I136").Value = 2268
$ws.Range("J136").Value = 41668896
$ws.Range("K136").Value = 6804
$ws.Range("L136").Value = 125006688
$ws.Range("M136").Value = -4254
$ws.Range("N136").Value = -125011788

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6472
$ws.Range("J62").Value = 6472
$ws.Range("L62").Value = 6472
$ws.Range("N62").Value = -7720
$ws.Range("H65").Value = 6472
$ws.Range("J65").Value = 6472
$ws.Range("L65").Value = 32360
$ws.Range("N65").Value = -38600
$ws.Range("H81").Value = 5374.222
$ws.Range("J81").Value = 4124.5
$ws.Range("L81").Value = 8249
$ws.Range("N81").Value = -10371
$ws.Range("H84").Value = 5374.222
$ws.Range("J84").Value = 4124.5
$ws.Range("L84").Value = 41245
$ws.Range("N84").Value = -51853
$ws.Range("H113").Value = 750.8889
$ws.Range("I113").Value = 694
$ws.Range("J113").Value = 950
$ws.Range("K113").Value = 2082
$ws.Range("L113").Value = 2850
$ws.Range("M113").Value = 88
$ws.Range("N113").Value = -7190
$ws.Range("H122").Value = 2800
$ws.Range("I122").Value = 2558.3333
$ws.Range("J122").Value = 3283.3333
$ws.Range("K122").Value = 7674.999899999999
$ws.Range("L122").Value = 9849.999899999999
$ws.Range("M122").Value = -5224.999899999999
$ws.Range("N122").Value = -14749.9999
